# Generate Report for Handback
# For each localized-language sheet (zh-cn, de-de), rows 2 and 3 represent
# the two source files that have now been handed back:
#  - Status (col B) moves from "Ready for handoff" to
#    "Handed back: in sync with en-US"
#  - Latest Target File (col E) is populated with the source file (same as
#    column A), as a hyperlink
#  - Latest Handback File (col F) is populated with the latest handoff
#    xlf file (same as column C), as a hyperlink
#  - Latest Handback DateTime (col G) is stamped with the handback time

$wb = $excel.ActiveWorkbook

function Update-HandbackSheet {
    param($SheetName, $HandbackDateTime)

    $ws = $wb.Worksheets.Item($SheetName)

    # Collect the existing hyperlink address/display text for the A and C
    # columns of rows 2 and 3 so the new E/F hyperlinks can reuse them.
    $linkInfo = @{}
    foreach ($hl in $ws.Hyperlinks) {
        $addr = $hl.Range.Address()
        $linkInfo[$addr] = @{ Url = $hl.Address; Display = $hl.TextToDisplay }
    }

    for ($row = 2; $row -le 3; $row++) {
        # Status -> Handed back: in sync with en-US
        $ws.Cells.Item($row, 2).Value = "Handed back: in sync with en-US"

        $aAddr = "`$A`$$row"
        $cAddr = "`$C`$$row"

        $aInfo = $linkInfo[$aAddr]
        $cInfo = $linkInfo[$cAddr]

        # Latest Target File (column E) - mirrors column A
        $eCell = $ws.Cells.Item($row, 5)
        $ws.Hyperlinks.Add($eCell, $aInfo.Url, [Type]::Missing, [Type]::Missing, $aInfo.Display) | Out-Null

        # Latest Handback File (column F) - mirrors column C
        $fCell = $ws.Cells.Item($row, 6)
        $ws.Hyperlinks.Add($fCell, $cInfo.Url, [Type]::Missing, [Type]::Missing, $cInfo.Display) | Out-Null

        # Latest Handback DateTime (column G)
        $ws.Cells.Item($row, 7).Value = $HandbackDateTime
    }
}

Update-HandbackSheet "zh-cn" "2016-03-08 12:08:39"
Update-HandbackSheet "de-de" "2016-03-08 12:08:49"
